$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove obsolete listings (rows 7-19) first, while row numbers are still "original" ---
$ws.Range("A7:H19").EntireRow.Delete()

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 42
$ws.Columns.Item(4).ColumnWidth = 30
$ws.Columns.Item(8).ColumnWidth = 12

# --- Refresh the 5 remaining listings with the newly scraped data ---

# Row 2
$ws.Range("A2").Value = "2026-01-30 06:44:14"
$ws.Range("B2").Value = "【Java経験者】4月開始/与野 官公庁向けマイグレーション案件"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5482097"
$ws.Range("G2").Value = 85
$ws.Range("H2").Value = "★Java"

# Row 3
$ws.Range("A3").Value = "2026-01-30 06:44:14"
$ws.Range("B3").Value = "WEBサーバーの管理、トラブル解決対応できる方を募集します!"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5481859"
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = "◇管理"

# Row 4
$ws.Range("A4").Value = "2026-01-30 06:44:14"
$ws.Range("B4").Value = "進行管理およびチームディレクションを担当"
$ws.Range("D4").Value = "~ 5,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = "◇管理"

# Row 5
$ws.Range("A5").Value = "2026-01-30 06:44:14"
$ws.Range("B5").Value = "【高スキル】Web3.0系プロダクトの上級エンジニア募集"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5481715"
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = ""

# Row 6
$ws.Range("A6").Value = "2026-01-30 06:44:14"
$ws.Range("B6").Value = "無人美容什器[ 決済 × IoT 連携の設計サポート ](※実装なし/スポット)"
$ws.Range("D6").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5481888"
$ws.Range("G6").Value = 18
$ws.Range("H6").Value = ""

# --- Hyperlinks: wipe the stale collection (old + now-orphaned rows 7-19 entries)
#     and recreate exactly the 5 that should remain, pointed at the refreshed URLs ---
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5482097")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5481859")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5418064")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5481715")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5481888")
